$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-07 Saturday" "2024-12-08 Sunday"
Replace-Text "422×3=1266" "817×9=7353"
Replace-Text "268×5=1340" "383×3=1149"
Replace-Text "644×9=5796" "145×7=1015"
Replace-Text "331×9=2979" "563×6=3378"
Replace-Text "191×3=573" "157×8=1256"
Replace-Text "325×3=975" "482×4=1928"
Replace-Text "608×4=2432" "404×5=2020"
Replace-Text "870×2=1740" "767×4=3068"
Replace-Text "649×2=1298" "601×3=1803"
Replace-Text "624×3=1872" "561×9=5049"
Replace-Text "243×8=1944" "856×3=2568"
Replace-Text "841×5=4205" "446×3=1338"
Replace-Text "609×6=3654" "662×7=4634"
Replace-Text "886×9=7974" "649×5=3245"
Replace-Text "831×6=4986" "348×3=1044"
Replace-Text "210×2=420" "575×3=1725"
Replace-Text "889×5=4445" "687×3=2061"
Replace-Text "786×4=3144" "720×5=3600"
Replace-Text "179×7=1253" "122×8=976"
Replace-Text "663×3=1989" "695×3=2085"
Replace-Text "166×6=996" "990×5=4950"
Replace-Text "599×8=4792" "420×5=2100"
Replace-Text "300×6=1800" "775×9=6975"
Replace-Text "520×4=2080" "480×5=2400"
Replace-Text "453×5=2265" "918×2=1836"

Write-Output "Done"
